$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (B5): "Lord Xu" -> "the Lord"
$ws.Range("B5").Value = "Besides, they’ve worked in the manor for years—the Lord would know them well."

# Row 7 (B7): "Lady Mei" -> "Madam Mei"
$ws.Range("B7").Value = "Madam Mei"

# Row 8 (B8): "Xu Ming" -> "Ming"
$ws.Range("B8").Value = "Ming"

# Clear the now-unused "Question-Meeting" tag cells in column E for rows 7 and 8
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()

# Update the active selection
$ws.Range("E10").Select() | Out-Null
